# Applies the "chore: update Sheets via scheduled runner" edit described by the
# XML diff: refreshed currentAveragePrice / LevePrice* / LeveProfit* figures
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 440.66666
$ws.Range("I28").Value = 482.875
$ws.Range("K28").Value = 482.875
$ws.Range("M28").Value = 2.125

# Row 32
$ws.Range("H32").Value = 3377.6
$ws.Range("J32").Value = 3377.6
$ws.Range("L32").Value = 3377.6
$ws.Range("N32").Value = -4029.6

# Row 33
$ws.Range("H33").Value = 8065199.5
$ws.Range("I33").Value = 10000340
$ws.Range("J33").Value = 2116.5
$ws.Range("K33").Value = 10000340
$ws.Range("L33").Value = 2116.5
$ws.Range("M33").Value = -10000111
$ws.Range("N33").Value = -2574.5

# Row 107
$ws.Range("H107").Value = 2298.25
$ws.Range("I107").Value = 691.6818
$ws.Range("K107").Value = 691.6818
$ws.Range("M107").Value = 1228.3182

# Row 111
$ws.Range("H111").Value = 2555
$ws.Range("I111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("M111").ClearContents()

# Row 116
$ws.Range("H116").Value = 3929.5
$ws.Range("I116").Value = 3628.4
$ws.Range("K116").Value = 3628.4
$ws.Range("M116").Value = -186.4000000000001

# Row 132
$ws.Range("H132").Value = 10782.421
$ws.Range("I132").Value = 3124.6
$ws.Range("K132").Value = 9373.799999999999
$ws.Range("M132").Value = -6843.799999999999

# Row 138
$ws.Range("H138").Value = 3976.6316
$ws.Range("I138").Value = 6539.5
$ws.Range("J138").Value = 3293.2
$ws.Range("K138").Value = 19618.5
$ws.Range("L138").Value = 9879.599999999999
$ws.Range("M138").Value = -14478.5
$ws.Range("N138").Value = -20159.6

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 847574.9
$ws.Range("I61").Value = 4403.8076
$ws.Range("J61").Value = 3283402.5
$ws.Range("K61").Value = 4403.8076
$ws.Range("L61").Value = 3283402.5
$ws.Range("M61").Value = -4191.8076
$ws.Range("N61").Value = -3283826.5

# Row 132
$ws.Range("H132").Value = 1081314.2
$ws.Range("I132").Value = 5755.08
$ws.Range("K132").Value = 17265.24
$ws.Range("M132").Value = -14735.24

# Row 136
$ws.Range("H136").Value = 847574.9
$ws.Range("I136").Value = 4403.8076
$ws.Range("J136").Value = 3283402.5
$ws.Range("K136").Value = 13211.4228
$ws.Range("L136").Value = 9850207.5
$ws.Range("M136").Value = -10661.4228
$ws.Range("N136").Value = -9855307.5

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 19916.816
$ws.Range("I20").Value = 6915.4
$ws.Range("J20").Value = 40445.367
$ws.Range("K20").Value = 6915.4
$ws.Range("L20").Value = 40445.367
$ws.Range("M20").Value = -6668.4
$ws.Range("N20").Value = -40939.367

# Row 99
$ws.Range("H99").Value = 39832.5
$ws.Range("J99").Value = 7248.75
$ws.Range("L99").Value = 7248.75
$ws.Range("N99").Value = -10244.75

# Row 134
$ws.Range("H134").Value = 21236.285
$ws.Range("I134").Value = 15094.909
$ws.Range("J134").Value = 43754.668
$ws.Range("K134").Value = 45284.727
$ws.Range("L134").Value = 131264.004
$ws.Range("M134").Value = -42749.727
$ws.Range("N134").Value = -136334.004

$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Range("H99").Value = 5416.8335
$ws.Range("I99").Value = 5416.8335
$ws.Range("K99").Value = 5416.8335
$ws.Range("M99").Value = -3918.8335

# Row 105
$ws.Range("H105").Value = 30452.5
$ws.Range("I105").Value = 38936.668
$ws.Range("J105").Value = 5000
$ws.Range("K105").Value = 38936.668
$ws.Range("L105").Value = 5000
$ws.Range("M105").Value = -37189.668
$ws.Range("N105").Value = -8494

# Row 107
$ws.Range("H107").Value = 1677.2727
$ws.Range("I107").Value = 1850
$ws.Range("K107").Value = 1850
$ws.Range("M107").Value = 70

# Row 126
$ws.Range("H126").Value = 5416.8335
$ws.Range("I126").Value = 5416.8335
$ws.Range("K126").Value = 16250.5005
$ws.Range("M126").Value = -13780.5005

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 24.5
$ws.Range("J12").Value = 34.5
$ws.Range("L12").Value = 103.5
$ws.Range("N12").Value = -449.5

# Row 24
$ws.Range("H24").Value = 911.6667
$ws.Range("J24").Value = 1166.6666
$ws.Range("L24").Value = 3499.9998
$ws.Range("N24").Value = -3959.9998

$ws = $wb.Worksheets.Item("GSM")
# Row 3
$ws.Range("H3").Value = 887563.3
$ws.Range("I3").Value = 219
$ws.Range("K3").Value = 219
$ws.Range("M3").Value = -103

# Row 97
$ws.Range("H97").Value = 2849.25
$ws.Range("I97").Value = 793.65
$ws.Range("K97").Value = 793.65
$ws.Range("M97").Value = -297.65

# Row 102
$ws.Range("H102").Value = 6593.7856
$ws.Range("I102").Value = 9337.375
$ws.Range("J102").Value = 2935.6667
$ws.Range("K102").Value = 9337.375
$ws.Range("L102").Value = 2935.6667
$ws.Range("M102").Value = -7715.375
$ws.Range("N102").Value = -6179.6667

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 800.75
$ws.Range("I22").Value = 401
$ws.Range("K22").Value = 401
$ws.Range("M22").Value = -106

# Row 27
$ws.Range("H27").Value = 800.75
$ws.Range("I27").Value = 401
$ws.Range("K27").Value = 401
$ws.Range("M27").Value = -294

# Row 40
$ws.Range("H40").Value = 3691.5557
$ws.Range("I40").Value = 3403.125
$ws.Range("J40").Value = 5999
$ws.Range("K40").Value = 3403.125
$ws.Range("L40").Value = 5999
$ws.Range("M40").Value = -3267.125
$ws.Range("N40").Value = -6271

# Row 61
$ws.Range("H61").Value = 3432.7222
$ws.Range("I61").Value = 2774.2144
$ws.Range("K61").Value = 2774.2144
$ws.Range("M61").Value = -2572.2144

# Row 94
$ws.Range("H94").Value = 1000000000
$ws.Range("J94").Value = 1000000000
$ws.Range("L94").Value = 1000000000
$ws.Range("N94").Value = -1000001352

# Row 100
$ws.Range("H100").Value = 3165
$ws.Range("I100").Value = 2976.7856
$ws.Range("J100").Value = 3541.4285
$ws.Range("K100").Value = 2976.7856
$ws.Range("L100").Value = 3541.4285
$ws.Range("M100").Value = -2435.7856
$ws.Range("N100").Value = -4623.4285

# Row 113
$ws.Range("H113").Value = 3432.7222
$ws.Range("I113").Value = 2774.2144
$ws.Range("K113").Value = 2774.2144
$ws.Range("M113").Value = -604.2143999999998

# Row 122
$ws.Range("H122").Value = 7626.8335
$ws.Range("I122").Value = 8164.7
$ws.Range("K122").Value = 24494.1
$ws.Range("M122").Value = -22044.1

# Row 132
$ws.Range("H132").Value = 1001024.4
$ws.Range("I132").Value = 4428.6
$ws.Range("K132").Value = 13285.8
$ws.Range("M132").Value = -10755.8

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 2500
$ws.Range("I122").Value = 2500
$ws.Range("K122").Value = 7500
$ws.Range("M122").Value = -5050
